# Update trajectory values on sheet "strategy_id-0" with data from Edmundo
# (Sunday Nov 3 2024) for gdp_mmm_usd, occrateinit_gnrl_occupancy,
# population_gnrl_rural and population_gnrl_urban (columns J:AS, i.e. the
# yearly trajectory values 0..35).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("strategy_id-0")

# Row 2 (gdp_mmm_usd): J2:AS2
$ws.Range("J2").Value = 883.625114435657
$ws.Range("K2").Value = 879.6633193356651
$ws.Range("L2").Value = 875.128803492218
$ws.Range("M2").Value = 869.209782818966
$ws.Range("N2").Value = 863.059317505741
$ws.Range("O2").Value = 865.689844847097
$ws.Range("P2").Value = 884.482880054306
$ws.Range("Q2").Value = 912.436872900765
$ws.Range("R2").Value = 940.915334498977
$ws.Range("S2").Value = 969.1250867557389
$ws.Range("T2").Value = 998.042563797904
$ws.Range("U2").Value = 1027.87933783827
$ws.Range("V2").Value = 1058.65576561112
$ws.Range("W2").Value = 1090.39220385073
$ws.Range("X2").Value = 1123.10900929138
$ws.Range("Y2").Value = 1156.82653866736
$ws.Range("Z2").Value = 1191.56514871294
$ws.Range("AA2").Value = 1227.34519616242
$ws.Range("AB2").Value = 1264.18703775006
$ws.Range("AC2").Value = 1302.11442415646
$ws.Range("AD2").Value = 1341.16468184739
$ws.Range("AE2").Value = 1381.37853123495
$ws.Range("AF2").Value = 1422.79669273123
$ws.Range("AG2").Value = 1465.4598867483
$ws.Range("AH2").Value = 1509.40883369826
$ws.Range("AI2").Value = 1554.6842539932
$ws.Range("AJ2").Value = 1601.3268680452
$ws.Range("AK2").Value = 1649.37739626635
$ws.Range("AL2").Value = 1698.87655906874
$ws.Range("AM2").Value = 1749.86507686445
$ws.Range("AN2").Value = 1802.38367006558
$ws.Range("AO2").Value = 1856.4730590842
$ws.Range("AP2").Value = 1912.17396433241
$ws.Range("AQ2").Value = 1969.5271062223
$ws.Range("AR2").Value = 2028.57320516594
$ws.Range("AS2").Value = 2089.35298157544

# Row 10 (occrateinit_gnrl_occupancy): J10:AS10
$ws.Range("J10").Value = 4.18121887287025
$ws.Range("K10").Value = 4.18121887287025
$ws.Range("L10").Value = 4.18121887287025
$ws.Range("M10").Value = 4.18121887287025
$ws.Range("N10").Value = 4.18121887287025
$ws.Range("O10").Value = 4.18121887287025
$ws.Range("P10").Value = 4.18121887287025
$ws.Range("Q10").Value = 4.18121887287025
$ws.Range("R10").Value = 4.18121887287025
$ws.Range("S10").Value = 4.18121887287025
$ws.Range("T10").Value = 4.18121887287025
$ws.Range("U10").Value = 4.18121887287025
$ws.Range("V10").Value = 4.18121887287025
$ws.Range("W10").Value = 4.18121887287025
$ws.Range("X10").Value = 4.18121887287025
$ws.Range("Y10").Value = 4.18121887287025
$ws.Range("Z10").Value = 4.18121887287025
$ws.Range("AA10").Value = 4.18121887287025
$ws.Range("AB10").Value = 4.18121887287025
$ws.Range("AC10").Value = 4.18121887287025
$ws.Range("AD10").Value = 4.18121887287025
$ws.Range("AE10").Value = 4.18121887287025
$ws.Range("AF10").Value = 4.18121887287025
$ws.Range("AG10").Value = 4.18121887287025
$ws.Range("AH10").Value = 4.18121887287025
$ws.Range("AI10").Value = 4.18121887287025
$ws.Range("AJ10").Value = 4.18121887287025
$ws.Range("AK10").Value = 4.18121887287025
$ws.Range("AL10").Value = 4.18121887287025
$ws.Range("AM10").Value = 4.18121887287025
$ws.Range("AN10").Value = 4.18121887287025
$ws.Range("AO10").Value = 4.18121887287025
$ws.Range("AP10").Value = 4.18121887287025
$ws.Range("AQ10").Value = 4.18121887287025
$ws.Range("AR10").Value = 4.18121887287025
$ws.Range("AS10").Value = 4.18121887287025

# Row 11 (population_gnrl_rural): J11:AS11
$ws.Range("J11").Value = 21860642.31736
$ws.Range("K11").Value = 21728030.6592
$ws.Range("L11").Value = 21621460.5824
$ws.Range("M11").Value = 21490943.75736
$ws.Range("N11").Value = 21316691.82552
$ws.Range("O11").Value = 21105851.10642
$ws.Range("P11").Value = 20756060.6466545
$ws.Range("Q11").Value = 20523016.9252806
$ws.Range("R11").Value = 20326508.4840811
$ws.Range("S11").Value = 20093211.2595196
$ws.Range("T11").Value = 19833681.14343
$ws.Range("U11").Value = 19516422.1989416
$ws.Range("V11").Value = 19178676.949651
$ws.Range("W11").Value = 18824075.5523725
$ws.Range("X11").Value = 18456236.6507488
$ws.Range("Y11").Value = 18078474.3291
$ws.Range("Z11").Value = 17668714.1604554
$ws.Range("AA11").Value = 17255197.6984047
$ws.Range("AB11").Value = 16837784.1409137
$ws.Range("AC11").Value = 16417134.2677346
$ws.Range("AD11").Value = 15994757.211225
$ws.Range("AE11").Value = 15549446.9391769
$ws.Range("AF11").Value = 15102124.470786
$ws.Range("AG11").Value = 14652465.2622843
$ws.Range("AH11").Value = 14200776.3422551
$ws.Range("AI11").Value = 13746888.5502
$ws.Range("AJ11").Value = 13275038.7289486
$ws.Range("AK11").Value = 12800625.8003402
$ws.Range("AL11").Value = 12323477.171434
$ws.Range("AM11").Value = 11843722.018246
$ws.Range("AN11").Value = 11361529.16901
$ws.Range("AO11").Value = 10865849.7149905
$ws.Range("AP11").Value = 10366618.0553635
$ws.Range("AQ11").Value = 9863852.524266081
$ws.Range("AR11").Value = 9357983.273673531
$ws.Range("AS11").Value = 8848548.10143001

# Row 12 (population_gnrl_urban): J12:AS12
$ws.Range("J12").Value = 60192665.68264
$ws.Range("K12").Value = 61457385.3408
$ws.Range("L12").Value = 62817579.4176
$ws.Range("M12").Value = 64123524.24264
$ws.Range("N12").Value = 65304836.17448
$ws.Range("O12").Value = 66375915.89358
$ws.Range("P12").Value = 67209058.3533455
$ws.Range("Q12").Value = 68422839.0747194
$ws.Range("R12").Value = 69776128.5159189
$ws.Range("S12").Value = 71021565.74048039
$ws.Range("T12").Value = 72187016.85657001
$ws.Range("U12").Value = 73298241.8010584
$ws.Range("V12").Value = 74343654.050349
$ws.Range("W12").Value = 75330884.4476275
$ws.Range("X12").Value = 76268850.3492512
$ws.Range("Y12").Value = 77166635.6709
$ws.Range("Z12").Value = 78059508.83954459
$ws.Range("AA12").Value = 78938035.3015953
$ws.Range("AB12").Value = 79799904.8590863
$ws.Range("AC12").Value = 80646479.7322654
$ws.Range("AD12").Value = 81483727.788775
$ws.Range("AE12").Value = 82330572.0608231
$ws.Range("AF12").Value = 83168507.52921399
$ws.Range("AG12").Value = 83994513.73771571
$ws.Range("AH12").Value = 84808883.6577449
$ws.Range("AI12").Value = 85609056.4498
$ws.Range("AJ12").Value = 86417571.2710515
$ws.Range("AK12").Value = 87214200.19965979
$ws.Range("AL12").Value = 87995848.828566
$ws.Range("AM12").Value = 88761136.981754
$ws.Range("AN12").Value = 89508796.83099
$ws.Range("AO12").Value = 90245050.2850095
$ws.Range("AP12").Value = 90961858.94463649
$ws.Range("AQ12").Value = 91655694.47573391
$ws.Range("AR12").Value = 92326327.7263265
$ws.Range("AS12").Value = 92964224.89857
